$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-79 down to 45-80
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record
$ws.Cells.Item(44, 1).Value = 5
$ws.Cells.Item(44, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value = "Maule"
$ws.Cells.Item(44, 4).Value = 44566
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 5).Value = 7
$ws.Cells.Item(44, 6).Value = 100112001
$ws.Cells.Item(44, 7).Value = "Berenjena"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 200
$ws.Cells.Item(44, 11).Value = 9000
$ws.Cells.Item(44, 12).Value = 9000
$ws.Cells.Item(44, 13).Value = 9000
$ws.Cells.Item(44, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(44, 15).Value = "Región del Maule"
$ws.Cells.Item(44, 16).Value = 180
$ws.Cells.Item(44, 17).Value = 50
$ws.Cells.Item(44, 18).Value = "Hortaliza"
